# edit.ps1 — reproduce the target commit against before.pptx
#
# The canonical OOXML diff this commit represents has two parts:
#
#   1. ppt/slides/slide16.xml: the table's <a:tableStyleId> changes from
#      {5245A7F7-F686-44A1-BFEC-C3A5153FDA34} ("Table_0", the style defined
#      in ppt/tableStyles.xml) to {032F6C3F-350F-4C05-AAE3-0C8F29A6F966}
#      (a built-in PowerPoint table style id not locally redefined).
#
#   2. ppt/theme/theme1.xml (the slide master's theme, "Integral") and
#      ppt/theme/theme2.xml (the notes master's theme, "Office Theme")
#      trade their 12 theme colours (dk1/lt1/dk2/lt2/accent1-6/hlink/
#      folHlink) — theme1 ends up with the colours theme2 used to have,
#      and vice versa. (fontScheme/fmtScheme and everything else in both
#      theme parts are byte-identical already, so this is the only real
#      content difference between the two theme parts.)
#
$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------

$targetSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $targetSlide.Shapes.Count; $i++) {
    $shp = $targetSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{032F6C3F-350F-4C05-AAE3-0C8F29A6F966}")
    }
}

# --- 2. Theme colours -------------------------------------------------------
# VBA-style RGB() packs bytes as 0x00BBGGRR, matching PowerPoint's
# ColorScheme.Colors(i).RGB convention.
function RGBFn($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

# The slide master's theme (theme1.xml, "Integral") picks up the colours
# that the notes master's theme ("Office Theme") used to have.
$masterColors = $p.SlideMaster.ColorScheme
$masterColors.Colors(1).RGB  = (RGBFn 0x00 0x00 0x00)   # dk1
$masterColors.Colors(2).RGB  = (RGBFn 0xFF 0xFF 0xFF)   # lt1
$masterColors.Colors(3).RGB  = (RGBFn 0x44 0x54 0x6A)   # dk2
$masterColors.Colors(4).RGB  = (RGBFn 0xE7 0xE6 0xE6)   # lt2
$masterColors.Colors(5).RGB  = (RGBFn 0x5B 0x9B 0xD5)   # accent1
$masterColors.Colors(6).RGB  = (RGBFn 0xED 0x7D 0x31)   # accent2
$masterColors.Colors(7).RGB  = (RGBFn 0xA5 0xA5 0xA5)   # accent3
$masterColors.Colors(8).RGB  = (RGBFn 0xFF 0xC0 0x00)   # accent4
$masterColors.Colors(9).RGB  = (RGBFn 0x44 0x72 0xC4)   # accent5
$masterColors.Colors(10).RGB = (RGBFn 0x70 0xAD 0x47)   # accent6
$masterColors.Colors(11).RGB = (RGBFn 0x05 0x63 0xC1)   # hlink
$masterColors.Colors(12).RGB = (RGBFn 0x95 0x4F 0x72)   # folHlink

# The notes master's theme (theme2.xml, "Office Theme") symmetrically picks
# up the colours the slide master's theme ("Integral") used to have.
$notesColors = $p.NotesMaster.ColorScheme
$notesColors.Colors(1).RGB  = (RGBFn 0x00 0x00 0x00)   # dk1
$notesColors.Colors(2).RGB  = (RGBFn 0xFF 0xFF 0xFF)   # lt1
$notesColors.Colors(3).RGB  = (RGBFn 0x45 0x5F 0x51)   # dk2
$notesColors.Colors(4).RGB  = (RGBFn 0xE3 0xDE 0xD1)   # lt2
$notesColors.Colors(5).RGB  = (RGBFn 0x99 0xCB 0x38)   # accent1
$notesColors.Colors(6).RGB  = (RGBFn 0x63 0xA5 0x37)   # accent2
$notesColors.Colors(7).RGB  = (RGBFn 0xE6 0xD0 0x24)   # accent3
$notesColors.Colors(8).RGB  = (RGBFn 0xCC 0x97 0x00)   # accent4
$notesColors.Colors(9).RGB  = (RGBFn 0x4E 0xB3 0xCF)   # accent5
$notesColors.Colors(10).RGB = (RGBFn 0x37 0x8D 0xA6)   # accent6
$notesColors.Colors(11).RGB = (RGBFn 0x6B 0x9F 0x25)   # hlink
$notesColors.Colors(12).RGB = (RGBFn 0xB2 0x6B 0x02)   # folHlink
